$wb = $excel.ActiveWorkbook

# Sheet ALC, row 18
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 500
$ws.Range("I18").Value = 500
$ws.Range("K18").Value = 500
$ws.Range("M18").Value = -216

# Sheet ARM, row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2033
$ws.Range("I61").Value = 1335.3182
$ws.Range("K61").Value = 1335.3182
$ws.Range("M61").Value = -1123.3182

# Sheet ARM, row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3329.0408
$ws.Range("I74").Value = 711.375
$ws.Range("J74").Value = 5842
$ws.Range("K74").Value = 711.375
$ws.Range("L74").Value = 5842
$ws.Range("M74").Value = 162.625
$ws.Range("N74").Value = -7590

# Sheet ARM, row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 3329.0408
$ws.Range("I77").Value = 711.375
$ws.Range("J77").Value = 5842
$ws.Range("K77").Value = 3556.875
$ws.Range("L77").Value = 29210
$ws.Range("M77").Value = 811.125
$ws.Range("N77").Value = -37946

# Sheet ARM, row 103
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H103").Value = 40000
$ws.Range("J103").Value = 40000
$ws.Range("L103").Value = 40000
$ws.Range("N103").Value = -42344

# Sheet ARM, row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2033
$ws.Range("I136").Value = 1335.3182
$ws.Range("K136").Value = 4005.9546
$ws.Range("M136").Value = -1455.9546

# Sheet BSM, row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2031.52
$ws.Range("I86").Value = 1846.1177
$ws.Range("J86").Value = 2425.5
$ws.Range("K86").Value = 1846.1177
$ws.Range("L86").Value = 2425.5
$ws.Range("M86").Value = -723.1177
$ws.Range("N86").Value = -4671.5

# Sheet BSM, row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2031.52
$ws.Range("I89").Value = 1846.1177
$ws.Range("J89").Value = 2425.5
$ws.Range("K89").Value = 9230.5885
$ws.Range("L89").Value = 12127.5
$ws.Range("M89").Value = -3614.5885
$ws.Range("N89").Value = -23359.5

# Sheet BSM, row 103
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H103").Value = 15497.75
$ws.Range("J103").Value = 15497.75
$ws.Range("L103").Value = 15497.75
$ws.Range("N103").Value = -17841.75

# Sheet CRP, row 6
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 927782.4
$ws.Range("I6").Value = 1133578.5
$ws.Range("J6").Value = 1700
$ws.Range("K6").Value = 1133578.5
$ws.Range("L6").Value = 1700
$ws.Range("M6").Value = -1133465.5
$ws.Range("N6").Value = -1926

# Sheet CRP, row 7
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 56.46154
$ws.Range("J7").Value = 57
$ws.Range("L7").Value = 57
$ws.Range("N7").Value = -283

# Sheet CRP, row 17
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 1300
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()

# Sheet CRP, row 25
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 1082.6666
$ws.Range("I25").Value = 1082.6666
$ws.Range("K25").Value = 1082.6666
$ws.Range("M25").Value = -908.6666

# Sheet CRP, row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1189.7587
$ws.Range("I31").Value = 1094.7916
$ws.Range("K31").Value = 1094.7916
$ws.Range("M31").Value = -799.7916

# Sheet CRP, row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1189.7587
$ws.Range("I34").Value = 1094.7916
$ws.Range("K34").Value = 1094.7916
$ws.Range("M34").Value = -892.7916

# Sheet CRP, row 41
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 11256
$ws.Range("I41").Value = 9500
$ws.Range("J41").Value = 12426.667
$ws.Range("K41").Value = 9500
$ws.Range("L41").Value = 12426.667
$ws.Range("M41").Value = -9072
$ws.Range("N41").Value = -13282.667

# Sheet CRP, row 50
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 6546
$ws.Range("J50").Value = 8092
$ws.Range("L50").Value = 8092
$ws.Range("N50").Value = -9342

# Sheet CRP, row 51
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 22480
$ws.Range("J51").Value = 22480
$ws.Range("L51").Value = 22480
$ws.Range("N51").Value = -23952

# Sheet CRP, row 59
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 11206.4
$ws.Range("I59").Value = 9055
$ws.Range("J59").Value = 11744.25
$ws.Range("K59").Value = 9055
$ws.Range("L59").Value = 11744.25
$ws.Range("M59").Value = -7910
$ws.Range("N59").Value = -14034.25

# Sheet CRP, row 60
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 6290.7
$ws.Range("J60").Value = 7388.7144
$ws.Range("L60").Value = 7388.7144
$ws.Range("N60").Value = -8410.714400000001

# Sheet CRP, row 61
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H61").Value = 22480
$ws.Range("J61").Value = 22480
$ws.Range("L61").Value = 22480
$ws.Range("N61").Value = -23176

# Sheet CUL, row 31
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H31").Value = 1433.3334
$ws.Range("I31").Value = 500
$ws.Range("J31").Value = 3300
$ws.Range("K31").Value = 1500
$ws.Range("L31").Value = 9900
$ws.Range("M31").Value = -1212
$ws.Range("N31").Value = -10476

# Sheet CUL, row 33
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 71428860
$ws.Range("I33").Value = 125000270
$ws.Range("J33").Value = 318
$ws.Range("K33").Value = 750001620
$ws.Range("L33").Value = 1908
$ws.Range("M33").Value = -750001337
$ws.Range("N33").Value = -2474

# Sheet CUL, row 35
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H35").Value = 1497.5
$ws.Range("J35").Value = 1497.5
$ws.Range("L35").Value = 4492.5
$ws.Range("N35").Value = -5068.5

# Sheet CUL, row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1138268.4
$ws.Range("J131").Value = 1371413.4
$ws.Range("L131").Value = 4114240.2
$ws.Range("N131").Value = -4124320.2

# Sheet CUL, row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 30303910
$ws.Range("J132").Value = 1643.1666
$ws.Range("L132").Value = 14788.4994
$ws.Range("N132").Value = -19848.4994

# Sheet GSM, row 68
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

# Sheet GSM, row 71
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

# Sheet GSM, row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 550003
$ws.Range("J80").Value = 550003
$ws.Range("L80").Value = 550003
$ws.Range("N80").Value = -551999

# Sheet GSM, row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 550003
$ws.Range("J83").Value = 550003
$ws.Range("L83").Value = 2750015
$ws.Range("N83").Value = -2759999

# Sheet GSM, row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2531.6562
$ws.Range("I132").Value = 2353.4285
$ws.Range("K132").Value = 7060.2855
$ws.Range("M132").Value = -4530.2855

# Sheet LTW, row 18
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 24005
$ws.Range("J18").Value = 24005
$ws.Range("L18").Value = 24005
$ws.Range("N18").Value = -24349

# Sheet LTW, row 20
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 2983.3333
$ws.Range("I20").Value = 3000
$ws.Range("J20").Value = 2900
$ws.Range("K20").Value = 3000
$ws.Range("L20").Value = 2900
$ws.Range("M20").Value = -2774
$ws.Range("N20").Value = -3352

# Sheet LTW, row 82
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1649.1428
$ws.Range("J82").Value = 1815
$ws.Range("L82").Value = 1815
$ws.Range("N82").Value = -2537

# Sheet LTW, row 85
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 1649.1428
$ws.Range("J85").Value = 1815
$ws.Range("L85").Value = 1815
$ws.Range("N85").Value = -4311

# Sheet LTW, row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 6072.2583
$ws.Range("I122").Value = 6381.591
$ws.Range("J122").Value = 5316.1113
$ws.Range("K122").Value = 19144.773
$ws.Range("L122").Value = 15948.3339
$ws.Range("M122").Value = -16694.773
$ws.Range("N122").Value = -20848.3339

# Sheet WVR, row 62
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 19999.875
$ws.Range("I62").Value = 15333.333
$ws.Range("J62").Value = 22799.8
$ws.Range("K62").Value = 15333.333
$ws.Range("L62").Value = 22799.8
$ws.Range("M62").Value = -14709.333
$ws.Range("N62").Value = -24047.8

# Sheet WVR, row 65
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 19999.875
$ws.Range("I65").Value = 15333.333
$ws.Range("J65").Value = 22799.8
$ws.Range("K65").Value = 76666.66500000001
$ws.Range("L65").Value = 113999
$ws.Range("M65").Value = -73546.66500000001
$ws.Range("N65").Value = -120239

# Sheet WVR, row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 618.3570999999999
$ws.Range("I107").Value = 358.23077
$ws.Range("J107").Value = 4000
$ws.Range("K107").Value = 1074.69231
$ws.Range("L107").Value = 12000
$ws.Range("M107").Value = 845.3076900000001
$ws.Range("N107").Value = -15840

# Sheet WVR, row 112
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H112").Value = 25088.666
$ws.Range("I112").Value = 20000
$ws.Range("J112").Value = 27633
$ws.Range("K112").Value = 20000
$ws.Range("L112").Value = 27633
$ws.Range("M112").Value = -18523
$ws.Range("N112").Value = -30587

# Sheet WVR, row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3197.0222
$ws.Range("I132").Value = 3647.6667
$ws.Range("J132").Value = 1957.75
$ws.Range("K132").Value = 10943.0001
$ws.Range("L132").Value = 5873.25
$ws.Range("M132").Value = -8413.000100000001
$ws.Range("N132").Value = -10933.25

# Sheet WVR, row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 10844.4
$ws.Range("I136").Value = 11257.263
$ws.Range("K136").Value = 33771.789
$ws.Range("M136").Value = -31221.789
